$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.75
$ws.Range("D2").Value = 0.8571428571428571
$ws.Range("E2").Value = 12
$ws.Range("B3").Value = 0.9
$ws.Range("C3").Value = 0.8181818181818182
$ws.Range("D3").Value = 0.8571428571428572
$ws.Range("E3").Value = 22
$ws.Range("B4").Value = 0.9375
$ws.Range("C4").Value = 0.9375
$ws.Range("D4").Value = 0.9375
$ws.Range("E4").Value = 16
$ws.Range("E5").Value = 12
$ws.Range("B6").Value = 0.9722222222222222
$ws.Range("D6").Value = 0.9859154929577464
$ws.Range("E6").Value = 35
$ws.Range("C7").Value = 0.6666666666666666
$ws.Range("D7").Value = 0.8
$ws.Range("E7").Value = 3
$ws.Range("C8").Value = 0.9166666666666666
$ws.Range("D8").Value = 0.9565217391304348
$ws.Range("E8").Value = 12
$ws.Range("B9").Value = 0.8148148148148148
$ws.Range("D9").Value = 0.8979591836734693
$ws.Range("E9").Value = 22
$ws.Range("C10").Value = 0.9090909090909091
$ws.Range("D10").Value = 0.9523809523809523
$ws.Range("E10").Value = 11
$ws.Range("B11").Value = 1
$ws.Range("C11").Value = 0.6666666666666666
$ws.Range("D11").Value = 0.8
$ws.Range("E11").Value = 6
$ws.Range("B12").Value = 0.9
$ws.Range("D12").Value = 0.9473684210526316
$ws.Range("E12").Value = 18
$ws.Range("E13").Value = 12
$ws.Range("C14").Value = 0.9166666666666666
$ws.Range("D14").Value = 0.9565217391304348
$ws.Range("E14").Value = 12
$ws.Range("E15").Value = 13
$ws.Range("C16").Value = 0.9285714285714286
$ws.Range("D16").Value = 0.962962962962963
$ws.Range("E16").Value = 14
$ws.Range("B17").Value = 0.3333333333333333
$ws.Range("C17").Value = 0.6666666666666666
$ws.Range("D17").Value = 0.4444444444444444
$ws.Range("E17").Value = 3
$ws.Range("C18").Value = 0.875
$ws.Range("D18").Value = 0.9333333333333333
$ws.Range("E18").Value = 8
$ws.Range("E19").Value = 5
$ws.Range("C20").Value = 0.2
$ws.Range("D20").Value = 0.3333333333333334
$ws.Range("E20").Value = 5
$ws.Range("E21").Value = 8
$ws.Range("E22").Value = 28
$ws.Range("B23").Value = 0.7058823529411765
$ws.Range("D23").Value = 0.8275862068965517
$ws.Range("E23").Value = 12
$ws.Range("B24").Value = 0.9090909090909091
$ws.Range("D24").Value = 0.9523809523809523
$ws.Range("E24").Value = 10
$ws.Range("B25").Value = 0.9297658862876255
$ws.Range("C25").Value = 0.9297658862876255
$ws.Range("D25").Value = 0.9297658862876255
$ws.Range("E25").Value = 0.9297658862876255
$ws.Range("B26").Value = 0.9336018970609764
$ws.Range("C26").Value = 0.8805077169207604
$ws.Range("D26").Value = 0.8870649772157809
$ws.Range("E26").Value = 299
$ws.Range("B27").Value = 0.9448668930766118
$ws.Range("C27").Value = 0.9297658862876255
$ws.Range("D27").Value = 0.9280713624955823
$ws.Range("E27").Value = 299
